# "Add files via upload" - turn the free-form lab report into a
# structured one: centered bold title, then Цели / Ход работы (itemised
# "-" bullet sentences, with QZipReader spell-check markers) / Итог
# sections. The trailing blank paragraph + sectPr are left untouched.
#
# The whole rewrite is expressed as a single block of WordprocessingML
# applied with Range.InsertXML, which is the only way to reproduce the
# exact run/proofErr/bookmark structure the document ends up with
# (Word's Range.Text setter only ever touches one run at a time here).

$d = $word.ActiveDocument

# Paragraphs 1-3 hold all of the original narrative text that gets
# redistributed below; paragraph 4 is the trailing blank paragraph that
# must stay untouched, so the replacement range stops right before it.
$originalContentEnd = $d.Paragraphs(3).Range.End
$target = $d.Range(0, $originalContentEnd)

$newBody = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:jc w:val="center"/><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="_Hlk120140881"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">Отчет по лабораторной работе </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>2</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Цели</w:t></w:r><w:r><w:t xml:space="preserve">: </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Получение навыков при работе с различными форматами хранения растровых изображений, получению информации об изображении, хранящемся в файле.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Ход работы</w:t></w:r><w:r><w:t xml:space="preserve">: </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>-Предварительно сделан макет интерфейса и функций.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>-</w:t></w:r><w:r><w:t>И</w:t></w:r><w:r><w:t>зучены основные форматы хранения растровых изображений и характеристики этих изображений.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>-</w:t></w:r><w:r><w:t xml:space="preserve">Изучены библиотеки для работы с файлами и изображениями, а также подключена отдельная библиотека </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>QZipReader</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>для работы с архивами.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>-Сделаны функции вы</w:t></w:r><w:r><w:t>вод</w:t></w:r><w:r><w:t>а</w:t></w:r><w:r><w:t xml:space="preserve"> врем</w:t></w:r><w:r><w:t>ени</w:t></w:r><w:r><w:t xml:space="preserve"> работы, </w:t></w:r><w:r><w:t>есть возможность</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>выб</w:t></w:r><w:r><w:t>о</w:t></w:r><w:r><w:t>ра</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> од</w:t></w:r><w:r><w:t>ного</w:t></w:r><w:r><w:t xml:space="preserve"> файл</w:t></w:r><w:r><w:t>а</w:t></w:r><w:r><w:t>, нескольк</w:t></w:r><w:r><w:t>их</w:t></w:r><w:r><w:t xml:space="preserve"> или архив</w:t></w:r><w:r><w:t>а</w:t></w:r><w:r><w:t xml:space="preserve"> с изображениями. </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>-</w:t></w:r><w:r><w:t>Сделан интерфейс для отображения всей информации в виде таблицы.</w:t></w:r><w:r><w:tab/></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>-</w:t></w:r><w:r><w:t xml:space="preserve">Сделан </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>exe</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">файл, а исходные файлы с документацией загружены на </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>git</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Итог</w:t></w:r><w:r><w:t>:</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Получены навыков при работе с различными форматами хранения растровых изображений, получению информации об изображении, хранящемся в файле. Изучены </w:t></w:r><w:r><w:t xml:space="preserve">способы подключения сторонних библиотек в </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Qt</w:t></w:r><w:r><w:t xml:space="preserve">. </w:t></w:r><w:r><w:t>Изучены основные характеристики изображений.</w:t></w:r><w:bookmarkStart w:id="1" w:name="_GoBack"/><w:bookmarkEnd w:id="1"/><w:bookmarkEnd w:id="0"/></w:p>'

$target.InsertXML($newBody)

Write-Output ("Paragraphs after rewrite: " + $d.Paragraphs.Count)
